# Scheduled-runner market data refresh: updates currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns (H:N) for a set of leve rows across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 123: market data is no longer available (price dropped to 0), so the
# profit figure in N123 can no longer be computed and is cleared entirely.
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H137").Value = 2870.742
$ws.Range("I137").Value = 1645.8572
$ws.Range("J137").Value = 5443
$ws.Range("K137").Value = 4937.571599999999
$ws.Range("L137").Value = 16329
$ws.Range("M137").Value = -2387.571599999999
$ws.Range("N137").Value = -21429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1477.1
$ws.Range("I2").Value = 1471.375
$ws.Range("K2").Value = 1471.375
$ws.Range("M2").Value = -1358.375

$ws.Range("H32").Value = 4641.2783
$ws.Range("I32").Value = 3687.7742
$ws.Range("K32").Value = 3687.7742
$ws.Range("M32").Value = -3400.7742

$ws.Range("H61").Value = 1465
$ws.Range("I61").Value = 1346.1111
$ws.Range("K61").Value = 1346.1111
$ws.Range("M61").Value = -1134.1111

$ws.Range("H74").Value = 4047.3845
$ws.Range("I74").Value = 3853.5652
$ws.Range("J74").Value = 5533.3335
$ws.Range("K74").Value = 3853.5652
$ws.Range("L74").Value = 5533.3335
$ws.Range("M74").Value = -2979.5652
$ws.Range("N74").Value = -7281.3335

$ws.Range("H77").Value = 4047.3845
$ws.Range("I77").Value = 3853.5652
$ws.Range("J77").Value = 5533.3335
$ws.Range("K77").Value = 19267.826
$ws.Range("L77").Value = 27666.6675
$ws.Range("M77").Value = -14899.826
$ws.Range("N77").Value = -36402.6675

$ws.Range("H116").Value = 1477.1
$ws.Range("I116").Value = 1471.375
$ws.Range("K116").Value = 1471.375
$ws.Range("M116").Value = 822.625

$ws.Range("H132").Value = 1786.1143
$ws.Range("I132").Value = 1043.8667
$ws.Range("J132").Value = 6239.6
$ws.Range("K132").Value = 3131.6001
$ws.Range("L132").Value = 18718.8
$ws.Range("M132").Value = -601.6001000000001
$ws.Range("N132").Value = -23778.8

$ws.Range("H136").Value = 1465
$ws.Range("I136").Value = 1346.1111
$ws.Range("K136").Value = 4038.3333
$ws.Range("M136").Value = -1488.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1477.1
$ws.Range("I3").Value = 1471.375
$ws.Range("K3").Value = 1471.375
$ws.Range("M3").Value = -1357.375

$ws.Range("H75").Value = 10038
$ws.Range("I75").Value = 6068.4
$ws.Range("J75").Value = 15000
$ws.Range("K75").Value = 6068.4
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -5132.4
$ws.Range("N75").Value = -16872

$ws.Range("H78").Value = 10038
$ws.Range("I78").Value = 6068.4
$ws.Range("J78").Value = 15000
$ws.Range("K78").Value = 18205.2
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -13525.2
$ws.Range("N78").Value = -54360

$ws.Range("H134").Value = 2017.6666
$ws.Range("I134").Value = 1410.1471
$ws.Range("J134").Value = 4599.625
$ws.Range("K134").Value = 4230.4413
$ws.Range("L134").Value = 13798.875
$ws.Range("M134").Value = -1695.4413
$ws.Range("N134").Value = -18868.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6173964
$ws.Range("I16").Value = 7937468.5
$ws.Range("K16").Value = 7937468.5
$ws.Range("M16").Value = -7937181.5

$ws.Range("H31").Value = 14288397
$ws.Range("I31").Value = 1272.9474
$ws.Range("J31").Value = 31254356
$ws.Range("K31").Value = 1272.9474
$ws.Range("L31").Value = 31254356
$ws.Range("M31").Value = -977.9474
$ws.Range("N31").Value = -31254946

$ws.Range("H34").Value = 14288397
$ws.Range("I34").Value = 1272.9474
$ws.Range("J34").Value = 31254356
$ws.Range("K34").Value = 1272.9474
$ws.Range("L34").Value = 31254356
$ws.Range("M34").Value = -1070.9474
$ws.Range("N34").Value = -31254760

$ws.Range("H58").Value = 1572.8202
$ws.Range("I58").Value = 1366.0122
$ws.Range("J58").Value = 3995.4285
$ws.Range("K58").Value = 1366.0122
$ws.Range("L58").Value = 3995.4285
$ws.Range("M58").Value = -1163.0122
$ws.Range("N58").Value = -4401.4285

$ws.Range("H113").Value = 6173964
$ws.Range("I113").Value = 7937468.5
$ws.Range("K113").Value = 7937468.5
$ws.Range("M113").Value = -7935298.5

$ws.Range("H122").Value = 2380
$ws.Range("I122").Value = 1118
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 3354
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -904
$ws.Range("N122").Value = -49900

$ws.Range("H132").Value = 2157.7693
$ws.Range("I132").Value = 1462.625
$ws.Range("J132").Value = 10499.5
$ws.Range("K132").Value = 4387.875
$ws.Range("L132").Value = 31498.5
$ws.Range("M132").Value = -1857.875
$ws.Range("N132").Value = -36558.5

$ws.Range("H134").Value = 6541.3184
$ws.Range("I134").Value = 8393.923000000001
$ws.Range("J134").Value = 3865.3333
$ws.Range("K134").Value = 25181.769
$ws.Range("L134").Value = 11595.9999
$ws.Range("M134").Value = -22646.769
$ws.Range("N134").Value = -16665.9999

$ws.Range("H136").Value = 1572.8202
$ws.Range("I136").Value = 1366.0122
$ws.Range("J136").Value = 3995.4285
$ws.Range("K136").Value = 4098.036599999999
$ws.Range("L136").Value = 11986.2855
$ws.Range("M136").Value = -1548.036599999999
$ws.Range("N136").Value = -17086.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8197589
$ws.Range("J131").Value = 937.94446
$ws.Range("L131").Value = 2813.83338
$ws.Range("N131").Value = -12893.83338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7215.9443
$ws.Range("I70").Value = 6007.25
$ws.Range("J70").Value = 9633.333000000001
$ws.Range("K70").Value = 6007.25
$ws.Range("L70").Value = 9633.333000000001
$ws.Range("M70").Value = -5737.25
$ws.Range("N70").Value = -10173.333

$ws.Range("H73").Value = 7215.9443
$ws.Range("I73").Value = 6007.25
$ws.Range("J73").Value = 9633.333000000001
$ws.Range("K73").Value = 6007.25
$ws.Range("L73").Value = 9633.333000000001
$ws.Range("M73").Value = -5071.25
$ws.Range("N73").Value = -11505.333

$ws.Range("H107").Value = 4831460
$ws.Range("I107").Value = 283.35294
$ws.Range("J107").Value = 18519792
$ws.Range("K107").Value = 283.35294
$ws.Range("L107").Value = 18519792
$ws.Range("M107").Value = 1636.64706
$ws.Range("N107").Value = -18523632

$ws.Range("H132").Value = 2938.95
$ws.Range("I132").Value = 1991.9333
$ws.Range("J132").Value = 5780
$ws.Range("K132").Value = 5975.7999
$ws.Range("L132").Value = 17340
$ws.Range("M132").Value = -3445.7999
$ws.Range("N132").Value = -22400

$ws.Range("H138").Value = 48480
$ws.Range("J138").Value = 48480
$ws.Range("L138").Value = 48480
$ws.Range("N138").Value = -58760

$ws.Range("H140").Value = 38750.77
$ws.Range("J140").Value = 38750.77
$ws.Range("L140").Value = 38750.77
$ws.Range("N140").Value = -49110.77

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4172.4287
$ws.Range("I7").Value = 2162.8
$ws.Range("J7").Value = 5288.8887
$ws.Range("K7").Value = 2162.8
$ws.Range("L7").Value = 5288.8887
$ws.Range("M7").Value = -2050.8
$ws.Range("N7").Value = -5512.8887

$ws.Range("H126").Value = 4172.4287
$ws.Range("I126").Value = 2162.8
$ws.Range("J126").Value = 5288.8887
$ws.Range("K126").Value = 6488.400000000001
$ws.Range("L126").Value = 15866.6661
$ws.Range("M126").Value = -4018.400000000001
$ws.Range("N126").Value = -20806.6661

$ws.Range("H132").Value = 12462.692
$ws.Range("I132").Value = 10961.2
$ws.Range("J132").Value = 50000
$ws.Range("K132").Value = 32883.60000000001
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -30353.60000000001
$ws.Range("N132").Value = -155060

$ws.Range("H136").Value = 4570.7144
$ws.Range("I136").Value = 2200
$ws.Range("J136").Value = 10497.5
$ws.Range("K136").Value = 6600
$ws.Range("L136").Value = 31492.5
$ws.Range("M136").Value = -4050
$ws.Range("N136").Value = -36592.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 309.85715
$ws.Range("I113").Value = 267.33334
$ws.Range("J113").Value = 366.55554
$ws.Range("K113").Value = 802.0000200000001
$ws.Range("L113").Value = 1099.66662
$ws.Range("M113").Value = 1367.99998
$ws.Range("N113").Value = -5439.66662

$ws.Range("H132").Value = 13890955
$ws.Range("I132").Value = 1241.2106
$ws.Range("J132").Value = 66671868
$ws.Range("K132").Value = 3723.6318
$ws.Range("L132").Value = 200015604
$ws.Range("M132").Value = -1193.6318
$ws.Range("N132").Value = -200020664

$ws.Range("H133").Value = 35166.92
$ws.Range("J133").Value = 35166.92
$ws.Range("L133").Value = 35166.92
$ws.Range("N133").Value = -45286.92

$ws.Range("H136").Value = 2467.25
$ws.Range("I136").Value = 1676.8572
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 5030.571599999999
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -2480.571599999999
$ws.Range("N136").Value = -29100
